# The author removed a post entry ("「自然界の消防士の訓練」...") that
# previously occupied row 505. Deleting the entire row shifts every
# subsequent row up by one (506->505, 507->506, ..., 580->579), which
# matches the diff exactly and also updates the sheet's used-range
# dimension from A1:C580 to A1:C579 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(505).Delete()
